$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Forestry, Trees & Timber Careers | AllAboutCareers"
$ws.Range("B2").Value = "https://www.allaboutcareers.com/careers/career-path/forestry-trees-timber"
$ws.Range("C2").Value = 538
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("A3").Value = "Visa"
$ws.Range("B3").Value = "https://www.tentsile.com/blogs/news/5-awesome-careers-for-people-who-love-trees?sa=X&ved=2ahUKEwjPhu3huZzmAhUNuZ4KHauQAZ4Q9QF6BAgLEAI"
$ws.Range("C3").Value = 56
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "jgpmc59@gmail.com`nzulujrmoses@gmail.com`ninfo@tentsile.com`nsupport@tentsile.com`nsupport@tentsile.com`nrepairs@tentsile.com`nrepairs@tentsile.com`ninfo@tentsile.com`n"

# Row 4
$ws.Range("A4").Value = "Tree Jobs, Employment | Indeed.com"
$ws.Range("B4").Value = "https://www.indeed.com/q-Tree-jobs.html"
$ws.Range("C4").Value = 1780
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("A5").Value = "Tree Service Jobs, Employment in Colorado | Indeed.com"
$ws.Range("B5").Value = "https://www.indeed.com/q-Tree-Service-l-Colorado-jobs.html"
$ws.Range("C5").Value = 1777
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("A6").Value = "50 Careers in Trees - Tree Foundation of Kern: About the Tree Foundation of Kern"
$ws.Range("B6").Value = "http://www.urbanforest.org/index.cfm/fuseaction/Pages.Page/id/430"
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("A7").Value = "Job Opportunities - Friends of Trees"
$ws.Range("B7").Value = "https://friendsoftrees.org/about/job-opportunities/"
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("A8").Value = "20 Best tree felling jobs (Hiring Now!) | SimplyHired"
$ws.Range("B8").Value = "https://www.simplyhired.com/search?q=tree+felling"
$ws.Range("C8").Value = 633
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 17
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("A9").Value = "Grist"
$ws.Range("B9").Value = "https://grist.org/article/2010-02-01-the-jobs-are-in-the-trees/"
$ws.Range("C9").Value = 74
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("A10").Value = "Green Jobs - Knowledge of Tree Care Career and Education."
$ws.Range("B10").Value = "https://californiareleaf.org/resources/green-jobs/"
$ws.Range("C10").Value = 64
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "lforlin@peralta.edu`nlforlin@peralta.edu`ncbrey@aplustree.com`ncbrey@aplustree.com`nckirkman@arborwell.com`nckirkman@arborwell.com`nandrew.misch@davey.com`nandrew.misch@davey.com`njbartolo@wcainc.com`njbartolo@wcainc.com`ncdiaz@wcainc.com`ncdiaz@wcainc.com`n"

# Delete row 11 entirely
$ws.Rows("11").Delete()

# Restore default row heights (avoid custom height caused by multi-line autofit)
$ws.Rows("3").AutoFit()
$ws.Rows("10").AutoFit()
